$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.129.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.303.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0820"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.22%  "

$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.655.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.312.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.805"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.059.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0918"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.42%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.51%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0740"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.106"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.27%  "

$ws.Range("E42").Value = "  -5.42%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.95%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.966.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("E48").Value = "  +20.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.79%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.525.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.39%  "
